$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Cells that only need a new number, style unchanged ---
$ws.Range("M14").Value = 33.333333333333
$ws.Range("G15").Value = 4
$ws.Range("N15").Value = -60.975609756097
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 86
$ws.Range("J16").Value = 91
$ws.Range("K16").Value = -5.494505494505
$ws.Range("L16").Value = -3.370786516853
$ws.Range("M16").Value = -23.893805309734
$ws.Range("N16").Value = -86.173633440514
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 71.428571428571
$ws.Range("I17").Value = 221
$ws.Range("J17").Value = 263
$ws.Range("K17").Value = -15.969581749049
$ws.Range("L17").Value = -19.636363636363
$ws.Range("M17").Value = 52.413793103448
$ws.Range("N17").Value = -53.571428571428
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 69
$ws.Range("J18").Value = 64
$ws.Range("K18").Value = 7.8125
$ws.Range("L18").Value = -9.210526315789
$ws.Range("M18").Value = -45.238095238095
$ws.Range("N18").Value = -90.322580645161
$ws.Range("C19").Value = 1
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = 157
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1.290322580645
$ws.Range("M19").Value = 70.652173913043
$ws.Range("N19").Value = -48.013245033112
$ws.Range("L20").Value = 15.217391304347
$ws.Range("N20").Value = -87.788018433179
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 40
$ws.Range("F21").Value = 53
$ws.Range("G21").Value = 41
$ws.Range("H21").Value = 29.268292682926
$ws.Range("I21").Value = 610
$ws.Range("J21").Value = 647
$ws.Range("K21").Value = -5.718701700154
$ws.Range("L21").Value = -8.408408408408
$ws.Range("M21").Value = 7.017543859649
$ws.Range("N21").Value = -76.574500768049
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 83
$ws.Range("K23").Value = -11.702127659574
$ws.Range("L23").Value = -2.352941176470
$ws.Range("M23").Value = 118.421052631579
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 42
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 533
$ws.Range("J24").Value = 556
$ws.Range("K24").Value = -4.136690647482
$ws.Range("L24").Value = -12.908496732026
$ws.Range("M24").Value = 61.515151515151
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -34.375
$ws.Range("I25").Value = 351
$ws.Range("J25").Value = 403
$ws.Range("K25").Value = -12.903225806451
$ws.Range("L25").Value = -6.896551724137
$ws.Range("M25").Value = -7.142857142857
$ws.Range("C26").Value = 3
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 32
$ws.Range("K26").Value = 10.344827586206
$ws.Range("L26").Value = -8.571428571428
$ws.Range("I27").Value = 47
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 6.818181818181
$ws.Range("L27").Value = 56.666666666666
$ws.Range("G28").Value = 1
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = -10.714285714285
$ws.Range("G29").Value = 1
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = -24

# --- Cells changing from text ("0"/"***.*", General style) to a number (style 15, #,##0) ---
# Use a same-style template cell (I14) to copy the number format without creating new styles.
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C18").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D27").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D28").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D29").Value = 1

# --- Cells changing from text ("0"/"***.*", General style) to a number (style 16, #,##0.0) ---
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E27").Value = 0
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E29").Value = -100

# --- Cells changing from a number (style 15) to text "0" (style 14) ---
# Copy formats then copy the whole cell (value+format) from C14, which already holds text "0" with style 14.
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# --- Cells changing from a number (style 16) to text "***.*" (style 14) ---
# Copy formats then copy the whole cell (value+format) from E14, which already holds text "***.*" with style 14.
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

